$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row above row 3 (shifts existing rows 3.. down by one)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new TreeID field definition
$ws.Cells.Item(3, 3).Value = "string"
$ws.Cells.Item(3, 5).Value = "TreeID"
$ws.Cells.Item(3, 6).Value = "TreeID"

# The "stand" field (now shifted to row 4) changes type from integer to string
$ws.Cells.Item(4, 3).Value = "string"

[void]$ws.Range("F3").Select()
